$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33
$ws.Cells.Item($row, 1).Value = 32
$ws.Cells.Item($row, 2).Value = "ProphetofTheChurch"
$ws.Cells.Item($row, 3).Value = "Prophet of The Church"
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = "Delusional Memories"
$ws.Cells.Item($row, 9).Value = 208
$ws.Cells.Item($row, 10).Value = 416
